$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6: update title and link
$ws.Range("D6").Value = "[Python pandas] Fill time interval by group :: 비어있는 시간 채우고, group별 직전값으로 NA 채우기"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Python-pandas-Fill-time-interval-by-group-%EB%B9%84%EC%96%B4%EC%9E%88%EB%8A%94-%EC%8B%9C%EA%B0%84-%EC%B1%84%EC%9A%B0%EA%B3%A0-group%EB%B3%84-%EC%A7%81%EC%A0%84%EA%B0%92%EC%9C%BC%EB%A1%9C-NA-%EC%B1%84%EC%9A%B0%EA%B8%B0"

# Row 28: update title and link
$ws.Range("D28").Value = "강화학습 : Offline Reinforcement Learning 자료 조사"
$ws.Range("E28").Value = "https://ropiens.tistory.com/142"

# Row 52: update title only
$ws.Range("D52").Value = "사건까지 걸린 시간은?"

$wb.Save()
